$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.21%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.66%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.031"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.37%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07896"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.17%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.843"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.82%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.106"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.01%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.788"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.15%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9195"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.21%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1343"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.06%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1897"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.55%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09108"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.52%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03468"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.12%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09836"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.11%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001414"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.25%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006145"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.10%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.722"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.56%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3439"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.05%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.00%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.167"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.63%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04405"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.62%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.68%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004614"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.58%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.86%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004446"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.10%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01940"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.29%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05079"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.63%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007592"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.64%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01020"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.51%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1344"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.95%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002153"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.40%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01019"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.11%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006178"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.43%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.85%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001661"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.37%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
